# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) onto the
# three new header cells so they match the other bold/centered headers.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Header labels for the new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-41) gets the team's season record.
$ws.Range("AD2:AD41").Value = 80
$ws.Range("AE2:AE41").Value = 82
$ws.Range("AF2:AF41").Value = 0
